# chore: update Sheets via scheduled runner
# Refresh market-price-derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) on several leve rows across all eight job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3789112.2
$ws.Range("I33").Value = 1412
$ws.Range("J33").Value = 12987813
$ws.Range("K33").Value = 1412
$ws.Range("L33").Value = 12987813
$ws.Range("M33").Value = -1183
$ws.Range("N33").Value = -12988271
$ws.Range("H140").Value = 85000
$ws.Range("J140").Value = 85000
$ws.Range("L140").Value = 85000
$ws.Range("N140").Value = -95360
$ws.Range("H141").Value = 1975.7179
$ws.Range("I141").Value = 1158.4
$ws.Range("J141").Value = 4700.1113
$ws.Range("K141").Value = 3475.2
$ws.Range("L141").Value = 14100.3339
$ws.Range("M141").Value = 1704.8
$ws.Range("N141").Value = -24460.3339

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 144022
$ws.Range("I2").Value = 182818.9
$ws.Range("J2").Value = 1766.6666
$ws.Range("K2").Value = 182818.9
$ws.Range("L2").Value = 1766.6666
$ws.Range("M2").Value = -182705.9
$ws.Range("N2").Value = -1992.6666
$ws.Range("H31").Value = 12596.917
$ws.Range("I31").Value = 7116.3
$ws.Range("J31").Value = 40000
$ws.Range("K31").Value = 7116.3
$ws.Range("L31").Value = 40000
$ws.Range("M31").Value = -6822.3
$ws.Range("N31").Value = -40588
$ws.Range("H32").Value = 9136.049999999999
$ws.Range("I32").Value = 7461.655
$ws.Range("J32").Value = 17926.625
$ws.Range("K32").Value = 7461.655
$ws.Range("L32").Value = 17926.625
$ws.Range("M32").Value = -7174.655
$ws.Range("N32").Value = -18500.625
$ws.Range("H74").Value = 1631.0294
$ws.Range("I74").Value = 1283.2354
$ws.Range("J74").Value = 1978.8235
$ws.Range("K74").Value = 1283.2354
$ws.Range("L74").Value = 1978.8235
$ws.Range("M74").Value = -409.2354
$ws.Range("N74").Value = -3726.8235
$ws.Range("H77").Value = 1631.0294
$ws.Range("I77").Value = 1283.2354
$ws.Range("J77").Value = 1978.8235
$ws.Range("K77").Value = 6416.177
$ws.Range("L77").Value = 9894.1175
$ws.Range("M77").Value = -2048.177
$ws.Range("N77").Value = -18630.1175
$ws.Range("H116").Value = 144022
$ws.Range("I116").Value = 182818.9
$ws.Range("J116").Value = 1766.6666
$ws.Range("K116").Value = 182818.9
$ws.Range("L116").Value = 1766.6666
$ws.Range("M116").Value = -180524.9
$ws.Range("N116").Value = -6354.6666
$ws.Range("H125").Value = 68583.25
$ws.Range("J125").Value = 68583.25
$ws.Range("L125").Value = 68583.25
$ws.Range("N125").Value = -78423.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 144022
$ws.Range("I3").Value = 182818.9
$ws.Range("J3").Value = 1766.6666
$ws.Range("K3").Value = 182818.9
$ws.Range("L3").Value = 1766.6666
$ws.Range("M3").Value = -182704.9
$ws.Range("N3").Value = -1994.6666
$ws.Range("H107").Value = 126854.29
$ws.Range("I107").Value = 178464.88
$ws.Range("J107").Value = 1514.2858
$ws.Range("K107").Value = 178464.88
$ws.Range("L107").Value = 1514.2858
$ws.Range("M107").Value = -176544.88
$ws.Range("N107").Value = -5354.2858
$ws.Range("H140").Value = 52089.168
$ws.Range("J140").Value = 52089.168
$ws.Range("L140").Value = 52089.168
$ws.Range("N140").Value = -62449.168

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3336228.2
$ws.Range("I6").Value = 5002500
$ws.Range("J6").Value = 3685
$ws.Range("K6").Value = 5002500
$ws.Range("L6").Value = 3685
$ws.Range("M6").Value = -5002387
$ws.Range("N6").Value = -3911
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -887
$ws.Range("N7").Value = -1226
$ws.Range("H25").Value = 1000000000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1000000000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1000000000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -1000000348
$ws.Range("H31").Value = 279301.44
$ws.Range("I31").Value = 1310.5385
$ws.Range("K31").Value = 1310.5385
$ws.Range("M31").Value = -1015.5385
$ws.Range("H34").Value = 279301.44
$ws.Range("I34").Value = 1310.5385
$ws.Range("K34").Value = 1310.5385
$ws.Range("M34").Value = -1108.5385
$ws.Range("H41").Value = 16200
$ws.Range("J41").Value = 19000
$ws.Range("L41").Value = 19000
$ws.Range("N41").Value = -19856
$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -26250
$ws.Range("H51").Value = 26500
$ws.Range("J51").Value = 26500
$ws.Range("L51").Value = 26500
$ws.Range("N51").Value = -27972
$ws.Range("H59").Value = 23629.875
$ws.Range("J59").Value = 23629.875
$ws.Range("L59").Value = 23629.875
$ws.Range("N59").Value = -25919.875
$ws.Range("H60").Value = 25000
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022
$ws.Range("H61").Value = 26500
$ws.Range("J61").Value = 26500
$ws.Range("L61").Value = 26500
$ws.Range("N61").Value = -27196
$ws.Range("H68").Value = 29332.666
$ws.Range("J68").Value = 29332.666
$ws.Range("L68").Value = 29332.666
$ws.Range("N68").Value = -30830.666
$ws.Range("H71").Value = 29332.666
$ws.Range("J71").Value = 29332.666
$ws.Range("L71").Value = 87997.99800000001
$ws.Range("N71").Value = -95485.99800000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2414.2856
$ws.Range("J39").Value = 2414.2856
$ws.Range("L39").Value = 7242.8568
$ws.Range("N39").Value = -7830.8568

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2351.279
$ws.Range("I102").Value = 2175.6775
$ws.Range("J102").Value = 2804.9167
$ws.Range("K102").Value = 2175.6775
$ws.Range("L102").Value = 2804.9167
$ws.Range("M102").Value = -553.6774999999998
$ws.Range("N102").Value = -6048.9167
$ws.Range("H126").Value = 11061.523
$ws.Range("I126").Value = 14286.134
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 42858.402
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -40388.402
$ws.Range("N126").Value = -13940
$ws.Range("H140").Value = 29597.166
$ws.Range("J140").Value = 29597.166
$ws.Range("L140").Value = 29597.166
$ws.Range("N140").Value = -39957.166

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5435018
$ws.Range("I122").Value = 7150827
$ws.Range("J122").Value = 2003399.8
$ws.Range("K122").Value = 21452481
$ws.Range("L122").Value = 6010199.4
$ws.Range("M122").Value = -21450031
$ws.Range("N122").Value = -6015099.4
$ws.Range("H136").Value = 11107.357
$ws.Range("I136").Value = 7435.087
$ws.Range("K136").Value = 22305.261
$ws.Range("M136").Value = -19755.261

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J113").Value = 3500
$ws.Range("L113").Value = 10500
$ws.Range("N113").Value = -14840
$ws.Range("H132").Value = 2108.919
$ws.Range("I132").Value = 1474.3077
$ws.Range("J132").Value = 3608.9092
$ws.Range("K132").Value = 4422.9231
$ws.Range("L132").Value = 10826.7276
$ws.Range("M132").Value = -1892.9231
$ws.Range("N132").Value = -15886.7276
$ws.Range("H136").Value = 2577.209
$ws.Range("I136").Value = 2446.5144
$ws.Range("K136").Value = 7339.5432
$ws.Range("M136").Value = -4789.5432
